# Apply updated dSF (column F) values per repulled data / mean calculation
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    5  = 2
    12 = -3
    13 = -1
    19 = 3
    27 = 2
    31 = 0
    33 = -1
    35 = 2
    39 = 0
    44 = 0
    51 = -1
    64 = -3
}

foreach ($row in $updates.Keys) {
    $ws.Cells.Item($row, 6).Value = $updates[$row]
}
